# Apply the edit described by the diff:
#  - HomePage!D1 gets a new numeric value (7), dimension/selection grow to include column D
#  - Users!C1 text changes from " Specialist" to " jamesthomas@gmail.com" (shared string content)
#  - Users column C widens to fit the new, longer text
#  - Users becomes the active sheet/tab (was ResetPage), selection becomes the whole column C
#  - ResetPage's sheetView loses tabSelected (handled automatically by switching the active sheet)

$wb = $excel.ActiveWorkbook

# --- HomePage: add a 4th column value and move the selection there ---
$wsHome = $wb.Worksheets.Item("HomePage")
$wsHome.Range("D1").Value = 7
$wsHome.Range("D1").Select()

# --- Users: update the label text, widen the column, select the whole column ---
$wsUsers = $wb.Worksheets.Item("Users")
$wsUsers.Range("C1").Value = " jamesthomas@gmail.com"
$wsUsers.Columns.Item(3).ColumnWidth = 23.8

# Make Users the active sheet/tab (moves tabSelected + bookViews.activeTab here)
$wsUsers.Activate()

# Select the entire column C on the now-active Users sheet
$wsUsers.Columns.Item(3).Select()
